# Automatische test-sync: 2025-06-24 21:18:50
# Adds a new log row (row 29) to the "Logs" sheet, extends the
# conditional-formatting ranges to include it, and bumps the
# "Retour / Terugbetaling" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$subject = "Maat ruilen mogelijk?"
$sender  = "mailmind.test@zohomail.eu"
$body    = "Hoi, ik heb per ongeluk de verkeerde maat ontvangen. Is het mogelijk om deze te ruilen voor de juiste maat?" + "`n" + "Sent using {0}"
$category = "Retour / Terugbetaling"
$reply   = "Beste klant," + "`n" + "Bedankt voor je bericht. Excuses voor het ongemak dat dit heeft veroorzaakt. Om je te helpen met het omruilen van de verkeerde maat voor de juiste maat, hebben we wat meer informatie nodig. Zou je alsjeblieft je ordernummer en de specifieke maat die je nodig hebt, kunnen doorgeven? Hiermee kunnen we de ruilprocedure voor je in gang zetten." + "`n" + "Verder, graag ontvangen we ook het artikel dat je hebt ontvangen in de verkeerde maat retour. We voorzien je van de instructies voor het retourneren zodra we bovenstaande informatie hebben ontvangen." + "`n" + "Bedankt voor je begrip en medewerking. We zien je reactie graag tegemoet." + "`n" + "Met vriendelijke groet," + "`n" + "[Naam bedrijf] E-mailassistent"
$timestamp = "2025-06-24 21:18:00"
$answered = "Ja"

$ws.Range("A29").Value = $subject
$ws.Range("B29").Value = $sender
$ws.Range("C29").Value = $body
$ws.Range("D29").Value = $category
$ws.Range("E29").Value = $reply
$ws.Range("F29").Value = $timestamp
$ws.Range("G29").Value = $answered

# Writing multi-line text bumps the row's height; restore it to match the
# sheet's standard (non-custom) row height, same as every other row.
$ws.Rows.Item(29).AutoFit()

# Extend the conditional-formatting ranges so they cover the new row too.
$fcsD = $ws.Range("D2:D28").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($ws.Range("D2:D29"))
}

$fcsG = $ws.Range("G2:G28").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($ws.Range("G2:G29"))
}

# Update the Dashboard summary count for "Retour / Terugbetaling".
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 12
